$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(2, @(2,9.490853973133028), @(3,7.059523881024897), @(4,7.230165426828693), @(5,11.73157900122326), @(6,38.72962792929074), @(8,7.344005520526261), @(9,30.23405911477384), @(11,10.65070691646894), @(13,14.67389947728568)),
  @(3, @(2,9.349053031241409), @(3,6.853847722206553), @(4,7.215854189163934), @(5,11.49052634643767), @(6,38.26644242949571), @(8,7.344005520526261), @(9,30.00185048550613), @(11,10.55337400510496), @(13,14.54821924202106)),
  @(4, @(2,9.265014795258494), @(3,6.726976198064682), @(4,7.206760901445862), @(5,11.34354466096585), @(6,37.98535806896166), @(8,7.344005520526261), @(9,29.86185222809561), @(11,10.49705343536182), @(13,14.47503259056664)),
  @(5, @(2,9.231584466157647), @(3,6.675220826125182), @(4,7.202978009991851), @(5,11.28399334106051), @(6,37.87174994336245), @(8,7.344005520526261), @(9,29.8054843989246), @(11,10.47499534734584), @(13,14.44623943025389)),
  @(6, @(2,9.226084221934849), @(3,6.666626419588972), @(4,7.202345185233942), @(5,11.27412827762267), @(6,37.85294489761799), @(8,7.344005520526261), @(9,29.79616676047193), @(11,10.47138731598166), @(13,14.44152144567734)),
  @(7, @(2,9.264560567851191), @(3,6.72627830018719), @(4,7.206710197493334), @(5,11.34274001899146), @(6,37.98382198597145), @(8,7.344005520526261), @(9,29.86108922411276), @(11,10.49675230181834), @(13,14.47464006373742)),
  @(8, @(2,9.441365741201418), @(3,6.988782632331001), @(4,7.225294254072629), @(5,11.64830309370899), @(6,38.56929616942895), @(8,7.344005520526261), @(9,30.15347924391639), @(11,10.6164491451988), @(13,14.62975710341379)),
  @(9, @(2,9.809637075753136), @(3,7.495071363077701), @(4,7.259323597300813), @(5,12.25188957793822), @(6,39.73910344683886), @(8,7.344005520526261), @(9,30.74576070239681), @(11,10.87718919878514), @(13,14.96412290620878)),
  @(10, @(2,10.08998293828777), @(3,7.857132262130587), @(4,7.282874863594794), @(5,12.69312999769564), @(6,40.60557248003014), @(8,7.344005520526261), @(9,31.19028275970065), @(11,11.08276397493688), @(13,15.22616885684334)),
  @(11, @(2,10.21892906082955), @(3,8.018810742887123), @(4,7.293278613833576), @(5,12.89240667693938), @(6,40.99995357179799), @(8,7.344005520526261), @(9,31.39406609418549), @(11,11.17890880643802), @(13,15.3484815013898)),
  @(12, @(2,10.26790286915725), @(3,8.079534345704808), @(4,7.297174014412657), @(5,12.96758360663901), @(6,41.14921750816891), @(8,7.344005520526261), @(9,31.47141773517905), @(11,11.21565811775076), @(13,15.39520507833234)),
  @(13, @(2,10.25734996710994), @(3,8.066479670605771), @(4,7.296337034881593), @(5,12.95140675591179), @(6,41.11707614151177), @(8,7.344005520526261), @(9,31.45475119735264), @(11,11.20772889969647), @(13,15.38512489131764)),
  @(14, @(2,10.22295559669147), @(3,8.023816895085718), @(4,7.293599976511578), @(5,12.89859766165071), @(6,41.0122358321031), @(8,7.344005520526261), @(9,31.40042642911888), @(11,11.18192556373004), @(13,15.35231756547202)),
  @(15, @(2,10.20190517778816), @(3,7.9976176896738), @(4,7.291917691693619), @(5,12.8662112650869), @(6,40.94800447544914), @(8,7.344005520526261), @(9,31.36717355539487), @(11,11.1661636388525), @(13,15.33227385492697)),
  @(16, @(2,10.08157942515319), @(3,7.846499579351908), @(4,7.282188743419779), @(5,12.68007103168743), @(6,40.57979362402677), @(8,7.344005520526261), @(9,31.17699299670829), @(11,11.07653044952934), @(13,15.21823452057557)),
  @(17, @(2,10.00808368230075), @(3,7.752970153796382), @(4,7.276141193124398), @(5,12.56545364517371), @(6,40.35388483872986), @(8,7.344005520526261), @(9,31.06069512301909), @(11,11.02219059195381), @(13,15.14904229988144)),
  @(18, @(2,9.965947538484672), @(3,7.698893309955588), @(4,7.272633622236022), @(5,12.49939700971851), @(6,40.22397535109284), @(8,7.344005520526261), @(9,30.99395380603571), @(11,10.99118530936737), @(13,15.10953963589323)),
  @(19, @(2,9.951706252334592), @(3,7.680537624500958), @(4,7.271440998675349), @(5,12.477011154117), @(6,40.17999833266003), @(8,7.344005520526261), @(9,30.97138340833535), @(11,10.98073142270708), @(13,15.09621657444386)),
  @(20, @(2,10.01589371634871), @(3,7.762956117759384), @(4,7.276787983150395), @(5,12.57766906765967), @(6,40.37793120399202), @(8,7.344005520526261), @(9,31.07305999525992), @(11,11.02794959682388), @(13,15.15637770186731)),
  @(21, @(2,10.23305459067286), @(3,8.036362062403251), @(4,7.29440511627312), @(5,12.91411729429611), @(6,41.0430330107204), @(8,7.344005520526261), @(9,31.41637829306769), @(11,11.18949565771952), @(13,15.36194316578643)),
  @(22, @(2,10.37579915618331), @(3,8.212105107874827), @(4,7.305660995430144), @(5,13.13231291977136), @(6,41.47720120374117), @(8,7.344005520526261), @(9,31.64180537875247), @(11,11.29704786123845), @(13,15.49864205475214)),
  @(23, @(2,10.29955791637638), @(3,8.118597138785216), @(4,7.299677030823132), @(5,13.01603697400838), @(6,41.24556077354937), @(8,7.344005520526261), @(9,31.52140875375349), @(11,11.23947688772008), @(13,15.42548170633239)),
  @(24, @(2,10.01236242934034), @(3,7.758442410078827), @(4,7.276495665183417), @(5,12.57214697257131), @(6,40.36705993185304), @(8,7.344005520526261), @(9,31.06746945983807), @(11,11.02534521462394), @(13,15.15306050118369)),
  @(25, @(2,9.708048747748261), @(3,7.359527348686844), @(4,7.250373940265181), @(5,12.08864971251935), @(6,39.42095730829759), @(8,7.344005520526261), @(9,30.58372714413319), @(11,10.80405287423356), @(13,14.87063865240734)),
)

foreach ($rowEntry in $data) {
    $r = $rowEntry[0]
    for ($i = 1; $i -lt $rowEntry.Count; $i++) {
        $pair = $rowEntry[$i]
        $c = $pair[0]
        $v = $pair[1]
        $ws.Cells.Item($r, $c).Value = $v
    }
}
